# Rename sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Horario"

# AutoFilter over the data range (creates autoFilter element)
$ws.Range("A1:P3").AutoFilter() | Out-Null

# Register the (hidden, sheet-scoped) _FilterDatabase defined name that Excel
# writes alongside an AutoFilter
$fd = $ws.Names.Add("_xlnm._FilterDatabase", "='Horario'!`$A`$1:`$P`$3")
$fd.Visible = $false

Write-Host "done"
